# Auto-generated PowerShell COM-interop script to apply scheduled market-data refresh
# to the Belias_Profits workbook (updates currentAveragePrice*/LevePrice*/LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 377.83334
$ws.Range("I33").Value = 188.875
$ws.Range("J33").Value = 1889.5
$ws.Range("K33").Value = 188.875
$ws.Range("L33").Value = 1889.5
$ws.Range("M33").Value = 40.125
$ws.Range("H74").Value = 3097.0908
$ws.Range("I74").Value = 2992.3076
$ws.Range("J74").Value = 3486.2856
$ws.Range("K74").Value = 2992.3076
$ws.Range("L74").Value = 3486.2856
$ws.Range("M74").Value = -2056.3076
$ws.Range("N74").Value = -5358.2856
$ws.Range("H76").Value = 3114.2856
$ws.Range("I76").Value = 3133.3333
$ws.Range("J76").Value = 3080
$ws.Range("K76").Value = 3133.3333
$ws.Range("L76").Value = 3080
$ws.Range("M76").Value = -2818.3333
$ws.Range("H77").Value = 3097.0908
$ws.Range("I77").Value = 2992.3076
$ws.Range("J77").Value = 3486.2856
$ws.Range("K77").Value = 14961.538
$ws.Range("L77").Value = 17431.428
$ws.Range("M77").Value = -10281.538
$ws.Range("N77").Value = -26791.428
$ws.Range("H79").Value = 3114.2856
$ws.Range("I79").Value = 3133.3333
$ws.Range("J79").Value = 3080
$ws.Range("K79").Value = 3133.3333
$ws.Range("L79").Value = 3080
$ws.Range("M79").Value = -2041.3333
$ws.Range("H100").Value = 2700.8333
$ws.Range("I100").Value = 2415.7144
$ws.Range("J100").Value = 3100
$ws.Range("K100").Value = 2415.7144
$ws.Range("L100").Value = 3100
$ws.Range("M100").Value = -1874.7144
$ws.Range("N100").Value = -4182
$ws.Range("H134").Value = 30944.875
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 30944.875
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 30944.875
$ws.Range("N134").Value = -41084.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6494594
$ws.Range("I45").Value = 9091831
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 9091831
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -9091454
$ws.Range("N45").Value = -2254
$ws.Range("H63").Value = 3463.5881
$ws.Range("I63").Value = 2440.7144
$ws.Range("J63").Value = 4179.6
$ws.Range("K63").Value = 2440.7144
$ws.Range("L63").Value = 4179.6
$ws.Range("M63").Value = -1754.7144
$ws.Range("N63").Value = -5551.6
$ws.Range("H66").Value = 3463.5881
$ws.Range("I66").Value = 2440.7144
$ws.Range("J66").Value = 4179.6
$ws.Range("K66").Value = 12203.572
$ws.Range("L66").Value = 20898
$ws.Range("M66").Value = -8771.572
$ws.Range("N66").Value = -27762

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 562.6667
$ws.Range("I22").Value = 474.8
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 474.8
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -301.8
$ws.Range("N22").Value = -1348
$ws.Range("H103").Value = 31250.25
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 31250.25
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 31250.25
$ws.Range("N103").Value = -33594.25
$ws.Range("H105").Value = 2364.4565
$ws.Range("I105").Value = 2358.2954
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2358.2954
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -611.2954
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3153.2354
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3153.2354
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3153.2354
$ws.Range("N4").Value = -3377.2354
$ws.Range("H62").Value = 2443.4375
$ws.Range("I62").Value = 2394.3333
$ws.Range("J62").Value = 2506.5715
$ws.Range("K62").Value = 2394.3333
$ws.Range("L62").Value = 2506.5715
$ws.Range("M62").Value = -1770.3333
$ws.Range("N62").Value = -3754.5715
$ws.Range("H65").Value = 2443.4375
$ws.Range("I65").Value = 2394.3333
$ws.Range("J65").Value = 2506.5715
$ws.Range("K65").Value = 11971.6665
$ws.Range("L65").Value = 12532.8575
$ws.Range("M65").Value = -8851.666499999999
$ws.Range("N65").Value = -18772.8575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17670.033
$ws.Range("I4").Value = 300.33334
$ws.Range("J4").Value = 19600
$ws.Range("K4").Value = 901.0000200000001
$ws.Range("L4").Value = 58800
$ws.Range("M4").Value = -789.0000200000001
$ws.Range("N4").Value = -59024
$ws.Range("H108").Value = 4200
$ws.Range("I108").Value = 3666.6667
$ws.Range("J108").Value = 5000
$ws.Range("K108").Value = 11000.0001
$ws.Range("L108").Value = 15000
$ws.Range("M108").Value = -8120.000100000001
$ws.Range("N108").Value = -20760
$ws.Range("H131").Value = 905.33
$ws.Range("I131").Value = 240
$ws.Range("J131").Value = 940.34735
$ws.Range("K131").Value = 720
$ws.Range("L131").Value = 2821.04205
$ws.Range("M131").Value = 4320
$ws.Range("N131").Value = -12901.04205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7725.5654
$ws.Range("I5").Value = 1568
$ws.Range("J5").Value = 8649.200000000001
$ws.Range("K5").Value = 1568
$ws.Range("L5").Value = 8649.200000000001
$ws.Range("M5").Value = -1456
$ws.Range("N5").Value = -8873.200000000001
$ws.Range("H70").Value = 5412.7646
$ws.Range("I70").Value = 5317.3335
$ws.Range("J70").Value = 5641.8
$ws.Range("K70").Value = 5317.3335
$ws.Range("L70").Value = 5641.8
$ws.Range("M70").Value = -5047.3335
$ws.Range("N70").Value = -6181.8
$ws.Range("H73").Value = 5412.7646
$ws.Range("I73").Value = 5317.3335
$ws.Range("J73").Value = 5641.8
$ws.Range("K73").Value = 5317.3335
$ws.Range("L73").Value = 5641.8
$ws.Range("M73").Value = -4381.3335
$ws.Range("N73").Value = -7513.8
$ws.Range("H80").Value = 2364.4
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 2222
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 2222
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -4218
$ws.Range("H83").Value = 2364.4
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 2222
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 11110
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -21094

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 100000
$ws.Range("I2").Value = 483333.34
$ws.Range("J2").Value = 11538.462
$ws.Range("K2").Value = 483333.34
$ws.Range("L2").Value = 11538.462
$ws.Range("M2").Value = -483221.34
$ws.Range("N2").Value = -11762.462
$ws.Range("H68").Value = 7391.1177
$ws.Range("I68").Value = 11074.9
$ws.Range("J68").Value = 2128.5715
$ws.Range("K68").Value = 11074.9
$ws.Range("L68").Value = 2128.5715
$ws.Range("M68").Value = -10325.9
$ws.Range("N68").Value = -3626.5715
$ws.Range("H71").Value = 7391.1177
$ws.Range("I71").Value = 11074.9
$ws.Range("J71").Value = 2128.5715
$ws.Range("K71").Value = 55374.5
$ws.Range("L71").Value = 10642.8575
$ws.Range("M71").Value = -51630.5
$ws.Range("N71").Value = -18130.8575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2985.7144
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 2983.3333
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 2983.3333
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -3207.3333
$ws.Range("H105").Value = 21833.334
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 21833.334
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 21833.334
$ws.Range("N105").Value = -28821.334
